{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Target change (per the supplied diff):\n//   Paragraph \"2.- Proceso de commit al repositorio.\" is split into runs so\n//   that \"commit\" is wrapped with <w:proofErr w:type=\"spellStart\"/> /\n//   <w:proofErr w:type=\"spellEnd\"/> (simulating Word's spellchecker flagging\n//   the English loanword), while the trailing \" \" run is left untouched.\n//\n//   Paragraph \"3.- Proceso de push.\" becomes\n//   '3.- Proceso de \"push\"' (curly quotes) with \"push\" similarly wrapped in\n//   proofErr spellStart/spellEnd markers, and the run with the old trailing\n//   period no longer exists (replaced by the closing curly quote run).\n//\n// We use Range.insertOoxml(..., Word.InsertLocation.replace) with a full\n// Flat-OPC package payload (the only format this host's insertOoxml\n// accepts) to splice in the exact run / proofErr structure, because plain\n// insertText() cannot create <w:proofErr/> markers or multiple runs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their (pre-edit) text rather than a\n// hard-coded index, so the script is resilient to unrelated paragraphs\n// being present before/after them.\nlet commitParagraph = null;\nlet pushParagraph = null;\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (t.indexOf(\"2.- Proceso de commit al repositorio.\") === 0) {\n    commitParagraph = p;\n  } else if (t.indexOf(\"3.- Proceso de push.\") === 0) {\n    pushParagraph = p;\n  }\n}\n\nif (!commitParagraph) {\n  throw new Error(\"Could not find paragraph '2.- Proceso de commit al repositorio.'\");\n}\nif (!pushParagraph) {\n  throw new Error(\"Could not find paragraph '3.- Proceso de push.'\");\n}\n\nconst FLAT_OPC_HEADER =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>';\nconst FLAT_OPC_FOOTER =\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n// --- Paragraph 2: \"2.- Proceso de commit al repositorio.\" -----------------\n// Keep the existing trailing \" \" run untouched by including it verbatim in\n// the replacement paragraph (insertOoxml Replace rewrites the whole range\n// it is called on, which here is the whole paragraph).\nconst commitParagraphOoxml =\n  FLAT_OPC_HEADER +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">2.- Proceso de </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>commit</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> al repositorio.</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"</w:p>\" +\n  FLAT_OPC_FOOTER;\n\ncommitParagraph.getRange().insertOoxml(commitParagraphOoxml, Word.InsertLocation.replace);\n\n// --- Paragraph 3: \"3.- Proceso de push.\" -----------------------------------\n// Becomes: 3.- Proceso de \"push\"   (curly quotes, no trailing period)\nconst pushParagraphOoxml =\n  FLAT_OPC_HEADER +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">3.- Proceso de </w:t></w:r>' +\n  \"<w:r><w:t>\\u201c</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>push</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"<w:r><w:t>\\u201d</w:t></w:r>\" +\n  \"</w:p>\" +\n  FLAT_OPC_FOOTER;\n\npushParagraph.getRange().insertOoxml(pushParagraphOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Target change (per the supplied diff):\n#   Paragraph \"2.- Proceso de commit al repositorio.\" is split into runs so\n#   that \"commit\" is wrapped with <w:proofErr w:type=\"spellStart\"/> /\n#   <w:proofErr w:type=\"spellEnd\"/> (simulating Word's spellchecker flagging\n#   the English loanword), while the trailing \" \" run is left untouched.\n#\n#   Paragraph \"3.- Proceso de push.\" becomes\n#   '3.- Proceso de \"push\"' (curly quotes) with \"push\" similarly wrapped in\n#   proofErr spellStart/spellEnd markers, and the run with the old trailing\n#   period no longer exists (replaced by the closing curly quote run).\n#\n# We use Range.InsertXML(...) with a full Flat-OPC package payload (the only\n# format this host's InsertXML accepts) to splice in the exact run /\n# proofErr structure, because plain text assignment cannot create\n# <w:proofErr/> markers or multiple runs.\n\n$d = $word.ActiveDocument\n\n$flatOpcHeader = '<?xml version=\"1.0\" standalone=\"yes\"?>' + `\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" ' + `\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n  '<pkg:xmlData>' + `\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n  '<w:body>'\n$flatOpcFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$openQuote = [char]0x201C\n$closeQuote = [char]0x201D\n\n# Locate the two target paragraphs by their (pre-edit) text rather than a\n# hard-coded index, so the script is resilient to unrelated paragraphs\n# being present before/after them.\n$commitParagraph = $null\n$pushParagraph = $null\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t.StartsWith(\"2.- Proceso de commit al repositorio.\")) {\n        $commitParagraph = $p\n    } elseif ($t.StartsWith(\"3.- Proceso de push.\")) {\n        $pushParagraph = $p\n    }\n}\n\nif ($commitParagraph -eq $null) {\n    throw \"Could not find paragraph '2.- Proceso de commit al repositorio.'\"\n}\nif ($pushParagraph -eq $null) {\n    throw \"Could not find paragraph '3.- Proceso de push.'\"\n}\n\n# --- Paragraph 2: \"2.- Proceso de commit al repositorio.\" ------------------\n# Keep the existing trailing \" \" run untouched by including it verbatim in\n# the replacement paragraph (InsertXML rewrites the whole range it is\n# called on, which here is the whole paragraph).\n$commitParagraphOoxml = $flatOpcHeader + `\n  '<w:p>' + `\n  '<w:r><w:t xml:space=\"preserve\">2.- Proceso de </w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellStart\"/>' + `\n  '<w:r><w:t>commit</w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellEnd\"/>' + `\n  '<w:r><w:t xml:space=\"preserve\"> al repositorio.</w:t></w:r>' + `\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n  '</w:p>' + `\n  $flatOpcFooter\n\n$commitParagraph.Range.InsertXML($commitParagraphOoxml)\n\n# --- Paragraph 3: \"3.- Proceso de push.\" ------------------------------------\n# Becomes: 3.- Proceso de \"push\"   (curly quotes, no trailing period)\n$pushParagraphOoxml = $flatOpcHeader + `\n  '<w:p>' + `\n  '<w:r><w:t xml:space=\"preserve\">3.- Proceso de </w:t></w:r>' + `\n  \"<w:r><w:t>$openQuote</w:t></w:r>\" + `\n  '<w:proofErr w:type=\"spellStart\"/>' + `\n  '<w:r><w:t>push</w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellEnd\"/>' + `\n  \"<w:r><w:t>$closeQuote</w:t></w:r>\" + `\n  '</w:p>' + `\n  $flatOpcFooter\n\n$pushParagraph.Range.InsertXML($pushParagraphOoxml)\n"}
